$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.863862156867981
$ws.Range("B1").Value = 1.858894467353821
$ws.Range("C1").Value = 1.931584596633911
$ws.Range("D1").Value = 3.224822998046875
$ws.Range("E1").Value = 4.314071655273438
